$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content fixes on the "2. Shop" module block and "Pets" module block.
# B11 was "brandName" - correct it to "name" (the shared string "brandName" is removed).
$ws.Range("B11").Value = "name"

# C26 was "userId" - correct it to "shopId" for the Owner field.
$ws.Range("C26").Value = "shopId"

# Scroll/selection state: the user selected the whole used range (e.g. Ctrl+A)
# and scrolled the view down so row 24 is at the top, leaving the prior active
# cell C26 as part of that selection.
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1:C37").Select()
